$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The manual review pass uppercased the lastname/firstname/middlename
# values in the first few data rows (rows 2-5).
$ws.Range("D2").Value = "ANDREWS"
$ws.Range("E2").Value = "CHANCE"
$ws.Range("D3").Value = "BANKS"
$ws.Range("E3").Value = "ANDREW"
$ws.Range("D4").Value = "BARTOLIN"
$ws.Range("E4").Value = "MATT"
$ws.Range("D5").Value = "BILLER"
$ws.Range("E5").Value = "VALERIE"
$ws.Range("F5").Value = "ANN"

# Update the last selected cell to match the reviewed workbook.
$ws.Range("D3").Select()
